# The source data gained one new weekly record. It is inserted as a new
# row 497 (pushing the previous rows 497-570 down to 498-571), then the
# new row's cells are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 497, shifting rows 497:570 down to 498:571.
$ws.Rows.Item(497).Insert()

# Populate the newly inserted row 497 with the new record's data.
$ws.Cells.Item(497, 1).Value = 4
$ws.Cells.Item(497, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(497, 3).Value = 'Los Lagos'
$ws.Cells.Item(497, 4).Value = 45077
$ws.Cells.Item(497, 5).Value = 10
$ws.Cells.Item(497, 6).Value = 100114013
$ws.Cells.Item(497, 7).Value = 'Zanahoria'
$ws.Cells.Item(497, 8).Value = 'Sin especificar'
$ws.Cells.Item(497, 9).Value = 'Primera'
$ws.Cells.Item(497, 10).Value = 70
$ws.Cells.Item(497, 11).Value = 8000
$ws.Cells.Item(497, 12).Value = 8000
$ws.Cells.Item(497, 13).Value = 8000
$ws.Cells.Item(497, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(497, 15).Value = 'Provincia de Llanquihue'
$ws.Cells.Item(497, 16).Value = 400
$ws.Cells.Item(497, 17).Value = 20
$ws.Cells.Item(497, 18).Value = 'Hortaliza'
